$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new header for column D
$ws.Cells.Item(1, 4).Value = "field_short_name"

# Map the existing "field" (column C) values to a shorter display name
# in the new "field_short_name" column (D), for every data row.
$lastRow = 29
for ($r = 2; $r -le $lastRow; $r++) {
    $field = $ws.Cells.Item($r, 3).Value()
    if ($field -eq "PX TO BOOK RATIO") {
        $short = "PB Ratio"
    } else {
        $short = "Price TR"
    }
    $ws.Cells.Item($r, 4).Value = $short
}

# Resize columns C and D to fit the new content
$ws.Columns.Item(3).ColumnWidth = 28.8
$ws.Columns.Item(4).ColumnWidth = 16.0

# Update the active selection to the new header cell
$ws.Range("D1").Select() | Out-Null
